$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark around "CMP73010" in the title line.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Replace the "Ben changing things up!" paragraph with the new text, and
#    drop the two empty trailing paragraphs that used to follow it, so the
#    new paragraph becomes the last paragraph of the body.
# ---------------------------------------------------------------------------

# Locate the paragraph that currently holds "Ben changing things up!" - it is
# the paragraph right after the ">>>  your stuff after this line >>>" one.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $paraText = $para.Range.Text.TrimEnd([char]13)
    if ($paraText -eq "Ben changing things up!") {
        $target = $para
        break
    }
}

if ($target -ne $null) {
    $prev = $target.Previous()

    # Delete everything from the end of the previous paragraph through to the
    # end of the document - this removes the old paragraph plus the two
    # trailing empty paragraphs (the very last paragraph mark of the body can
    # never be deleted, so this leaves a single, empty final paragraph).
    $wholeEnd = $d.Content.End
    $killRange = $d.Range($prev.Range.End, $wholeEnd)
    $killRange.Delete()

    # Add a fresh paragraph after the remaining last paragraph to host the
    # new content.
    $lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $endRange = $lastPara.Range
    $endRange.Collapse(0)
    $endRange.InsertParagraphAfter()

    # Populate the new, final paragraph with the exact run/formatting
    # structure via OOXML insertion (keeps the per-run language/font marks).
    $newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $newRange = $newPara.Range
    $xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:lang w:eastAsia="zh-CN"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">This </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="eastAsia"/>
      <w:lang w:eastAsia="zh-CN"/>
    </w:rPr>
    <w:t>is</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:eastAsia="zh-CN"/>
    </w:rPr>
    <w:t xml:space="preserve"> ysun38 Assignment 1 for Managing Software Development.</w:t>
  </w:r>
</w:p>
"@
    $newRange.InsertXML($xml)
}
